$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18 (item G=5471) on ALC
$ws.Cells.Item(18, 8).Value = 1953.4615
$ws.Cells.Item(18, 9).Value = 486.875
$ws.Cells.Item(18, 10).Value = 4300
$ws.Cells.Item(18, 11).Value = 486.875
$ws.Cells.Item(18, 12).Value = 4300
$ws.Cells.Item(18, 13).Value = -202.875
$ws.Cells.Item(18, 14).Value = -4868

# Row 32 (item G=5484) on ALC
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(32, 9).Value = 0
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 0
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 13).ClearContents()
$ws.Cells.Item(32, 14).ClearContents()

# Row 40 (item G=5505) on ALC
$ws.Cells.Item(40, 8).Value = 1947.8649
$ws.Cells.Item(40, 9).Value = 1902.84
$ws.Cells.Item(40, 10).Value = 2041.6666
$ws.Cells.Item(40, 11).Value = 1902.84
$ws.Cells.Item(40, 12).Value = 2041.6666
$ws.Cells.Item(40, 13).Value = -1727.84
$ws.Cells.Item(40, 14).Value = -2391.6666

# Row 116 (item G=27778) on ALC
$ws.Cells.Item(116, 8).Value = 2086044.5
$ws.Cells.Item(116, 9).Value = 2383628.2
$ws.Cells.Item(116, 10).Value = 2958
$ws.Cells.Item(116, 11).Value = 2383628.2
$ws.Cells.Item(116, 12).Value = 2958
$ws.Cells.Item(116, 13).Value = -2380186.2
$ws.Cells.Item(116, 14).Value = -9842

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (item G=44147) on ARM
$ws.Cells.Item(32, 8).Value = 3746.86
$ws.Cells.Item(32, 9).Value = 3365.848
$ws.Cells.Item(32, 10).Value = 8128.5
$ws.Cells.Item(32, 11).Value = 3365.848
$ws.Cells.Item(32, 12).Value = 8128.5
$ws.Cells.Item(32, 13).Value = -3078.848
$ws.Cells.Item(32, 14).Value = -8702.5

# Row 97 (item G=19941) on ARM
$ws.Cells.Item(97, 8).Value = 2783.125
$ws.Cells.Item(97, 9).Value = 3435.8333
$ws.Cells.Item(97, 10).Value = 825
$ws.Cells.Item(97, 11).Value = 3435.8333
$ws.Cells.Item(97, 12).Value = 825
$ws.Cells.Item(97, 13).Value = -2939.8333
$ws.Cells.Item(97, 14).Value = -1817

$ws = $wb.Worksheets.Item("BSM")
# Row 86 (item G=12526) on BSM
$ws.Cells.Item(86, 8).Value = 1898.4688
$ws.Cells.Item(86, 9).Value = 1824.1482
$ws.Cells.Item(86, 10).Value = 2299.8
$ws.Cells.Item(86, 11).Value = 1824.1482
$ws.Cells.Item(86, 12).Value = 2299.8
$ws.Cells.Item(86, 13).Value = -701.1482000000001
$ws.Cells.Item(86, 14).Value = -4545.8

# Row 89 (item G=12526) on BSM
$ws.Cells.Item(89, 8).Value = 1898.4688
$ws.Cells.Item(89, 9).Value = 1824.1482
$ws.Cells.Item(89, 10).Value = 2299.8
$ws.Cells.Item(89, 11).Value = 9120.741
$ws.Cells.Item(89, 12).Value = 11499
$ws.Cells.Item(89, 13).Value = -3504.741
$ws.Cells.Item(89, 14).Value = -22731

# Row 134 (item G=43998) on BSM
$ws.Cells.Item(134, 8).Value = 401817.84
$ws.Cells.Item(134, 9).Value = 667643.0600000001
$ws.Cells.Item(134, 10).Value = 3080
$ws.Cells.Item(134, 11).Value = 2002929.18
$ws.Cells.Item(134, 12).Value = 9240
$ws.Cells.Item(134, 13).Value = -2000394.18
$ws.Cells.Item(134, 14).Value = -14310

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (item G=5361) on CRP
$ws.Cells.Item(7, 8).Value = 47.0625
$ws.Cells.Item(7, 9).Value = 30.4
$ws.Cells.Item(7, 10).Value = 74.833336
$ws.Cells.Item(7, 11).Value = 30.4
$ws.Cells.Item(7, 12).Value = 74.833336
$ws.Cells.Item(7, 13).Value = 82.59999999999999
$ws.Cells.Item(7, 14).Value = -300.833336

# Row 32 (item G=2246) on CRP
$ws.Cells.Item(32, 8).Value = 1515.5714
$ws.Cells.Item(32, 9).Value = 1515.5714
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 1515.5714
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 13).Value = -1199.5714
$ws.Cells.Item(32, 14).ClearContents()

# Row 47 (item G=1920) on CRP
$ws.Cells.Item(47, 8).Value = 19999.5
$ws.Cells.Item(47, 9).Value = 0
$ws.Cells.Item(47, 10).Value = 19999.5
$ws.Cells.Item(47, 11).Value = 0
$ws.Cells.Item(47, 12).Value = 19999.5
$ws.Cells.Item(47, 13).ClearContents()
$ws.Cells.Item(47, 14).Value = -21131.5

# Row 132 (item G=44019) on CRP
$ws.Cells.Item(132, 8).Value = 1793.1714
$ws.Cells.Item(132, 9).Value = 1625.6538
$ws.Cells.Item(132, 10).Value = 2277.111
$ws.Cells.Item(132, 11).Value = 4876.9614
$ws.Cells.Item(132, 12).Value = 6831.333
$ws.Cells.Item(132, 13).Value = -2346.9614
$ws.Cells.Item(132, 14).Value = -11891.333

$ws = $wb.Worksheets.Item("CUL")
# Row 2 (item G=4847) on CUL
$ws.Cells.Item(2, 8).Value = 99024.55
$ws.Cells.Item(2, 9).Value = 165016.75
$ws.Cells.Item(2, 10).Value = 36.25
$ws.Cells.Item(2, 11).Value = 990100.5
$ws.Cells.Item(2, 12).Value = 217.5
$ws.Cells.Item(2, 13).Value = -989987.5
$ws.Cells.Item(2, 14).Value = -443.5

# Row 12 (item G=4854) on CUL
$ws.Cells.Item(12, 8).Value = 47.615383
$ws.Cells.Item(12, 9).Value = 72.125
$ws.Cells.Item(12, 10).Value = 36.72222
$ws.Cells.Item(12, 11).Value = 216.375
$ws.Cells.Item(12, 12).Value = 110.16666
$ws.Cells.Item(12, 13).Value = -43.375
$ws.Cells.Item(12, 14).Value = -456.16666

# Row 36 (item G=4732) on CUL
$ws.Cells.Item(36, 8).Value = 164.4
$ws.Cells.Item(36, 9).Value = 164.4
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 493.2
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).Value = -324.2
$ws.Cells.Item(36, 14).ClearContents()

# Row 39 (item G=4712) on CUL
$ws.Cells.Item(39, 8).Value = 3266
$ws.Cells.Item(39, 9).Value = 933.3333
$ws.Cells.Item(39, 10).Value = 3849.1667
$ws.Cells.Item(39, 11).Value = 2799.9999
$ws.Cells.Item(39, 12).Value = 11547.5001
$ws.Cells.Item(39, 13).Value = -2505.9999
$ws.Cells.Item(39, 14).Value = -12135.5001

# Row 44 (item G=4702) on CUL
$ws.Cells.Item(44, 8).Value = 3422.7
$ws.Cells.Item(44, 9).Value = 1461
$ws.Cells.Item(44, 10).Value = 8000
$ws.Cells.Item(44, 11).Value = 4383
$ws.Cells.Item(44, 12).Value = 24000
$ws.Cells.Item(44, 13).Value = -3985
$ws.Cells.Item(44, 14).Value = -24796

# Row 114 (item G=27865) on CUL
$ws.Cells.Item(114, 8).Value = 939.0333000000001
$ws.Cells.Item(114, 9).Value = 1600.1818
$ws.Cells.Item(114, 10).Value = 556.2632
$ws.Cells.Item(114, 11).Value = 4800.5454
$ws.Cells.Item(114, 12).Value = 1668.7896
$ws.Cells.Item(114, 13).Value = -1546.5454
$ws.Cells.Item(114, 14).Value = -8176.7896

# Row 117 (item G=27870) on CUL
$ws.Cells.Item(117, 8).Value = 2374.65
$ws.Cells.Item(117, 9).Value = 259.6
$ws.Cells.Item(117, 10).Value = 3079.6667
$ws.Cells.Item(117, 11).Value = 778.8000000000001
$ws.Cells.Item(117, 12).Value = 9239.000100000001
$ws.Cells.Item(117, 13).Value = 2663.2
$ws.Cells.Item(117, 14).Value = -16123.0001

# Row 121 (item G=27878) on CUL
$ws.Cells.Item(121, 8).Value = 1105.5714
$ws.Cells.Item(121, 9).Value = 339.82352
$ws.Cells.Item(121, 10).Value = 1512.375
$ws.Cells.Item(121, 11).Value = 1019.47056
$ws.Cells.Item(121, 12).Value = 4537.125
$ws.Cells.Item(121, 13).Value = 290.52944
$ws.Cells.Item(121, 14).Value = -7157.125

# Row 122 (item G=36078) on CUL
$ws.Cells.Item(122, 8).Value = 1304.2439
$ws.Cells.Item(122, 9).Value = 364.5909
$ws.Cells.Item(122, 10).Value = 2392.2632
$ws.Cells.Item(122, 11).Value = 3281.3181
$ws.Cells.Item(122, 12).Value = 21530.3688
$ws.Cells.Item(122, 13).Value = -831.3181
$ws.Cells.Item(122, 14).Value = -26430.3688

# Row 132 (item G=43972) on CUL
$ws.Cells.Item(132, 8).Value = 1185.1052
$ws.Cells.Item(132, 9).Value = 886.8570999999999
$ws.Cells.Item(132, 10).Value = 1359.0834
$ws.Cells.Item(132, 11).Value = 7981.7139
$ws.Cells.Item(132, 12).Value = 12231.7506
$ws.Cells.Item(132, 13).Value = -5451.7139
$ws.Cells.Item(132, 14).Value = -17291.7506

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (item G=14146) on GSM
$ws.Cells.Item(70, 8).Value = 2105733
$ws.Cells.Item(70, 9).Value = 5885310.5
$ws.Cells.Item(70, 10).Value = 5967.6665
$ws.Cells.Item(70, 11).Value = 5885310.5
$ws.Cells.Item(70, 12).Value = 5967.6665
$ws.Cells.Item(70, 13).Value = -5885040.5
$ws.Cells.Item(70, 14).Value = -6507.6665

# Row 73 (item G=14146) on GSM
$ws.Cells.Item(73, 8).Value = 2105733
$ws.Cells.Item(73, 9).Value = 5885310.5
$ws.Cells.Item(73, 10).Value = 5967.6665
$ws.Cells.Item(73, 11).Value = 5885310.5
$ws.Cells.Item(73, 12).Value = 5967.6665
$ws.Cells.Item(73, 13).Value = -5884374.5
$ws.Cells.Item(73, 14).Value = -7839.6665

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (item G=36249) on LTW
$ws.Cells.Item(7, 8).Value = 1837.3334
$ws.Cells.Item(7, 9).Value = 1358.3636
$ws.Cells.Item(7, 10).Value = 2590
$ws.Cells.Item(7, 11).Value = 1358.3636
$ws.Cells.Item(7, 12).Value = 2590
$ws.Cells.Item(7, 13).Value = -1246.3636
$ws.Cells.Item(7, 14).Value = -2814

# Row 16 (item G=5289) on LTW
$ws.Cells.Item(16, 8).Value = 860
$ws.Cells.Item(16, 9).Value = 850
$ws.Cells.Item(16, 10).Value = 866.6667
$ws.Cells.Item(16, 11).Value = 850
$ws.Cells.Item(16, 12).Value = 866.6667
$ws.Cells.Item(16, 13).Value = -680
$ws.Cells.Item(16, 14).Value = -1206.6667

# Row 38 (item G=2767) on LTW
$ws.Cells.Item(38, 8).Value = 10388.667
$ws.Cells.Item(38, 9).Value = 7000
$ws.Cells.Item(38, 10).Value = 12083
$ws.Cells.Item(38, 11).Value = 7000
$ws.Cells.Item(38, 12).Value = 12083
$ws.Cells.Item(38, 13).Value = -6590
$ws.Cells.Item(38, 14).Value = -12903

# Row 68 (item G=12563) on LTW
$ws.Cells.Item(68, 8).Value = 265404.78
$ws.Cells.Item(68, 9).Value = 715229.9
$ws.Cells.Item(68, 10).Value = 3006.8333
$ws.Cells.Item(68, 11).Value = 715229.9
$ws.Cells.Item(68, 12).Value = 3006.8333
$ws.Cells.Item(68, 13).Value = -714480.9
$ws.Cells.Item(68, 14).Value = -4504.8333

# Row 71 (item G=12563) on LTW
$ws.Cells.Item(71, 8).Value = 265404.78
$ws.Cells.Item(71, 9).Value = 715229.9
$ws.Cells.Item(71, 10).Value = 3006.8333
$ws.Cells.Item(71, 11).Value = 3576149.5
$ws.Cells.Item(71, 12).Value = 15034.1665
$ws.Cells.Item(71, 13).Value = -3572405.5
$ws.Cells.Item(71, 14).Value = -22522.1665

# Row 126 (item G=36249) on LTW
$ws.Cells.Item(126, 8).Value = 1837.3334
$ws.Cells.Item(126, 9).Value = 1358.3636
$ws.Cells.Item(126, 10).Value = 2590
$ws.Cells.Item(126, 11).Value = 4075.0908
$ws.Cells.Item(126, 12).Value = 7770
$ws.Cells.Item(126, 13).Value = -1605.0908
$ws.Cells.Item(126, 14).Value = -12710

$ws = $wb.Worksheets.Item("WVR")
# Row 122 (item G=36208) on WVR
$ws.Cells.Item(122, 8).Value = 1265.3125
$ws.Cells.Item(122, 9).Value = 1267.3077
$ws.Cells.Item(122, 10).Value = 1256.6666
$ws.Cells.Item(122, 11).Value = 3801.9231
$ws.Cells.Item(122, 12).Value = 3769.9998
$ws.Cells.Item(122, 13).Value = -1351.9231
$ws.Cells.Item(122, 14).Value = -8669.9998

# Row 126 (item G=36210) on WVR
$ws.Cells.Item(126, 8).Value = 3614.0908
$ws.Cells.Item(126, 9).Value = 4556.875
$ws.Cells.Item(126, 10).Value = 1100
$ws.Cells.Item(126, 11).Value = 13670.625
$ws.Cells.Item(126, 12).Value = 3300
$ws.Cells.Item(126, 13).Value = -11200.625
$ws.Cells.Item(126, 14).Value = -8240

Write-Output "Applied all updates"